# Update a few data points in row 2 of the training schedule and move the
# active selection to the cell that was last edited (E2), mirroring the
# manual edits made in Excel for this task run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 7
$ws.Range("G2").Value = 3
$ws.Range("H2").Value = 13

# Leave the selection on the last-touched cell, as the author did.
$ws.Range("E2").Select()
